$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$old1 = "✅ 1000 Bs = 6.0 = 24601.08 pesos"
$new1 = "✅ 1000 Bs = 6.01 = 24311.91 pesos"
$old2 = "✅ 24601.08 pesos = 5.97 = 966.62 Bs"
$new2 = "✅ 24311.91 pesos = 5.96 = 970.84 Bs"

$text = [string]$wsHoja1.Range("A1").Value2
$text = $text.Replace($old1, $new1)
$text = $text.Replace($old2, $new2)
$wsHoja1.Range("A1").Value = $text

# --- Update the rate figures on the "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 166.3
$wsTasas.Range("O10").Value = 4043.07
$wsTasas.Range("N12").Value = 4076.99
$wsTasas.Range("O12").Value = 162.805
